$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated computed values (cost_3 data refresh)
$ws.Range("A1").Value = 161.52538238402889
$ws.Range("B1").Value = 4.9610049959474267
$ws.Range("C1").Value = 3.9955968688845402

# Column C widened to match columns A/B (11.7109375 characters)
$ws.Columns.Item(3).ColumnWidth = 10.8
